$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 669 (shifts existing row 669 and everything
# below it down by one, turning the old last row 725 into row 726).
$ws.Rows(669).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A669").Value = 3
$ws.Range("B669").Value = "Femacal de La Calera"
$ws.Range("C669").Value = "Coquimbo"
$ws.Range("D669").Value = 45013
$ws.Range("E669").Value = 5
$ws.Range("F669").Value = 100112003
$ws.Range("G669").Value = "Ajo"
$ws.Range("H669").Value = "Chino"
$ws.Range("I669").Value = "Primera"
$ws.Range("J669").Value = 78
$ws.Range("K669").Value = 14500
$ws.Range("L669").Value = 15000
$ws.Range("M669").Value = 14744
$ws.Range("N669").Value = "`$/caja 10 kilos"
$ws.Range("O669").Value = "China"
$ws.Range("P669").Value = 1474
$ws.Range("Q669").Value = 10
$ws.Range("R669").Value = "Hortaliza"
